# Commit: Commit after including new fixed effects and running analysis as
# per potential responses to Maluccio comments
#
# 1. On the "growth" sheet, insert a new row at row 66 for the variable
#    v714a ("Maternal has a job but absent"), pushing the existing rows
#    (old 66-81) down to 67-82.
# 2. On the "v024 comparison" sheet, clear the AutoFilter criteria that was
#    filtering column F ("release_phase") down to "Phase 2" only, which
#    un-hides all the rows that filter had hidden.
# 3. Make "growth" the active sheet/tab again (it was "v024 comparison").

$wb = $excel.ActiveWorkbook

# --- 1. growth sheet: insert new variable row -----------------------------
$growth = $wb.Worksheets.Item("growth")

$growth.Rows.Item(66).Insert()

$growth.Range("A66").Value = "Maternal has a job but absent"
$growth.Range("B66").Value = "v714a"
$growth.Range("D66").Value = "v714a"

# --- 2. v024 comparison sheet: drop the "Phase 2" filter on column F ------
$v024 = $wb.Worksheets.Item("v024 comparison")

# Field 6 is column F (release_phase); toggling it off removes the
# <filterColumn>/<filters> criteria while keeping the AutoFilter range.
$v024.Range("A1:F38").AutoFilter(6)

# Removing the filter criteria should unhide the previously filtered rows;
# make sure every row is visible either way.
$v024.Cells.EntireRow.Hidden = $False

# --- 3. Restore "growth" as the active sheet / selection ------------------
$growth.Activate()
$growth.Range("C45").Select()
